$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "Tên phân loại" (category id) column L
# for rows 101-199 with their corresponding category values.
$ws.Range("L101:L106").Value = 6
$ws.Range("L107:L119").Value = 5
$ws.Range("L120:L144").Value = 7
$ws.Range("L145:L159").Value = 8
$ws.Range("L160:L180").Value = 9
$ws.Range("L181:L195").Value = 10
$ws.Range("L196:L199").Value = 11

# Row 22's wrap-text height was recalculated (autofit) after the edits above.
$ws.Rows.Item(22).RowHeight = 30.6

# Restore the scroll position (drop the stale topLeftCell) and leave the
# active selection where the user's last edit was made.
$ws.Range("N199").Select()
